$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2026-01-01 / Campbell / Holiday entry entirely (invoice feature
# removed -> that line item no longer belongs in the weekly export). Deleting
# the whole row shifts every row below it up by one.
$ws.Range("A2:F2").EntireRow.Delete()

# The 2026-01-02 entry (now row 2) changes client to Ueltschi with reduced
# hours/rate/total; date and type stay the same.
$ws.Range("B2").Value = "Ueltschi"
$ws.Range("C2").Value = 6.5
$ws.Range("E2").Value = 85
$ws.Range("F2").Value = 552.5

# SUBTOTAL row (now row 4) reflects the reduced hours/total.
$ws.Range("C4").Value = 6.5
$ws.Range("D4").Value = "Reg: 6.5 / OT: 0"
$ws.Range("F4").Value = 552.5
